# Regenerate orders with updated distance/size codes.
#
# The stimulus labels encode a trial's viewing distance and face size,
# e.g. "Face14_D80_S25", "Face14_D80_S25_l.png", "D80", "S30".
# This edit renumbers the distance codes and bumps one size code:
#   D80 -> D86
#   D64 -> D69
#   D51 -> D55
#   S30 -> S31
# (S20 and S25 are unchanged.) The substitutions are applied to every
# textual cell in the sheet, since the same codes recur across the
# Condition / Filename_Left / Filename_Right / Distance / Size columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$firstRow = $used.Row
$firstCol = $used.Column

$changed = 0

for ($r = $firstRow; $r -lt ($firstRow + $rowCount); $r++) {
    for ($c = $firstCol; $c -lt ($firstCol + $colCount); $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -is [string]) {
            $nv = $v.Replace("D80", "D86").Replace("D64", "D69").Replace("D51", "D55").Replace("S30", "S31")
            if ($nv -ne $v) {
                $cell.Value = $nv
                $changed = $changed + 1
            }
        }
    }
}

Write-Output "Updated $changed string cell(s) with new distance/size codes."
